# "Generate Report for Handoff"
#
# The file a4fc9b30-0d1e-495e-a482-ca587c3d247f.md was just handed off again
# (its handback version was stale), so its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", its handoff
# timestamp is refreshed, and an error detail describing the stale handback
# file is recorded - on both locale detail sheets (zh-cn, de-de) as well as
# being reflected back up on the Overview rollup sheet.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55792540f2923c0fcc3108a0d38a783322b26c2b/e2e/a4fc9b30-0d1e-495e-a482-ca587c3d247f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/168dea03647351622f05f4b22186ca443b677f30/e2e/a4fc9b30-0d1e-495e-a482-ca587c3d247f.md."

# ---- Overview sheet: row 3 is the a4fc9b30... file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = "2016-08-15 16:45:24"

# ---- zh-cn detail sheet: row 3 is the a4fc9b30... file ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = "2016-08-15 16:45:20"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de detail sheet: row 3 is the a4fc9b30... file ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = "2016-08-15 16:45:24"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.17
